$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Insert a new row before row 13 (old rows 13-28 shift down to 14-29)
$ws.Rows("13:13").Insert()

# Seed the new row's formatting from the row above (border / font / wrap all match
# the rest of the data table), then special-case a couple of cells below.
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's values (SEQ 5 - "Enable" flag column)
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Enable"
$ws.Range("C13").Value = "啟用記號"
$ws.Range("D13").Value = "VARCHAR2"
$ws.Range("E13").Value = 1
$ws.Range("G13").Value = "Y:啟用" + [char]10 + "N:停用"

# B13 ("Enable") is highlighted with a yellow fill
$ws.Range("B13").Interior.Color = 65535

# D13 ("VARCHAR2") is left aligned / vertically centered, not wrapped
$ws.Range("D13").HorizontalAlignment = -4131
$ws.Range("D13").VerticalAlignment = -4108
$ws.Range("D13").WrapText = $false

# E13 keeps the plain column look (centered, no wrap) instead of the table row style
$ws.Range("E13").HorizontalAlignment = -4108
$ws.Range("E13").VerticalAlignment = -4108
$ws.Range("E13").WrapText = $false
$ws.Range("E13").Font.Name = "標楷體"

# Row grows to fit the two-line note in G13
$ws.Rows("13:13").RowHeight = 32.4

# DBD becomes the active sheet/tab, with A18 selected
$ws.Activate()
$ws.Range("A18").Select()
